$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "44.593.46"
$ws.Range("E2").Value = "  +3.64%  "

$ws.Range("D3").Value = "2.423.82"
$ws.Range("E3").Value = "  +2.43%  "

$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  +0.03%  "

$ws.Range("D5").Value = "'312.83"
$ws.Range("E5").Value = "  +3.41%  "

$ws.Range("D6").Value = "'101.55"
$ws.Range("E6").Value = "  +5.71%  "

$ws.Range("E7").Value = "  +1.70%  "

$ws.Range("E8").Value = "  -0.03%  "

$ws.Range("D9").Value = "'0.516"
$ws.Range("E9").Value = "  +5.11%  "

$ws.Range("D10").Value = "'35.25"
$ws.Range("E10").Value = "  +3.14%  "

$ws.Range("B11").Value = "TRON"
$ws.Range("C11").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D11").Value = "'0.126"
$ws.Range("E11").Value = "  +1.87%  "

$ws.Range("B12").Value = "Dogecoin"
$ws.Range("C12").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D12").Value = "'0.0801"
$ws.Range("E12").Value = "  +1.97%  "

$ws.Range("D13").Value = "'18.98"
$ws.Range("E13").Value = "  +3.25%  "

$ws.Range("E14").Value = "  +2.73%  "

$ws.Range("D15").Value = "2.802.70"
$ws.Range("E15").Value = "  +2.61%  "

$ws.Range("D16").Value = "2.433.85"
$ws.Range("E16").Value = "  +3.80%  "

$ws.Range("E17").Value = "  +4.79%  "

$ws.Range("D18").Value = "44.479.82"
$ws.Range("E18").Value = "  +3.48%  "

$ws.Range("D19").Value = "'12.52"
$ws.Range("E19").Value = "  +4.52%  "

$ws.Range("D20").Value = "'6.39"
$ws.Range("E20").Value = "  +1.97%  "

$ws.Range("D21").Value = "0.0₃0929"
$ws.Range("E21").Value = "  +4.80%  "

$ws.Range("D22").Value = "'68.85"
$ws.Range("E22").Value = "  +1.32%  "

$ws.Range("D23").Value = "'241.18"
$ws.Range("E23").Value = "  +2.66%  "

$ws.Range("E24").Value = "  +4.63%  "

$ws.Range("E25").Value = "  +1.46%  "

$ws.Range("E26").Value = "  -0.07%  "

$ws.Range("E27").Value = "  +2.50%  "

$ws.Range("E28").Value = "  -4.20%  "

$ws.Range("D29").Value = "'9.60"
$ws.Range("E29").Value = "  +3.46%  "

$ws.Range("D30").Value = "'33.40"
$ws.Range("E30").Value = "  +5.04%  "

$ws.Range("D31").Value = "'48.64"
$ws.Range("E31").Value = "  +1.42%  "

$ws.Range("E32").Value = "  +18.03%  "

$ws.Range("D33").Value = "'19.56"
$ws.Range("E33").Value = "  +12.63%  "

$ws.Range("E34").Value = "  +2.98%  "

$ws.Range("E35").Value = "  +0.27%  "

$ws.Range("D36").Value = "'0.0769"
$ws.Range("E36").Value = "  +6.99%  "

$ws.Range("D37").Value = "'4.53"
$ws.Range("E37").Value = "  +3.66%  "

$ws.Range("E38").Value = "  +2.41%  "

$ws.Range("D39").Value = "'2.89"
$ws.Range("E39").Value = "  +3.65%  "

$ws.Range("D40").Value = "'124.04"
$ws.Range("E40").Value = "  +0.67%  "

$ws.Range("E41").Value = "  +1.04%  "

$ws.Range("D42").Value = "'2.20"
$ws.Range("E42").Value = "  -4.14%  "

$ws.Range("D43").Value = "'21.23"
$ws.Range("E43").Value = "  -1.93%  "

$ws.Range("E44").Value = "  +3.77%  "

$ws.Range("D45").Value = "1.949.99"
$ws.Range("E45").Value = "  +0.61%  "

$ws.Range("B46").Value = "ApeXProtocol"
$ws.Range("C46").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D46").Value = "'2.18"
$ws.Range("E46").Value = "  +2.17%  "

$ws.Range("B47").Value = "NEARProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D47").Value = "'2.93"
$ws.Range("E47").Value = "  +7.33%  "

$ws.Range("B48").Value = "FraxShare"
$ws.Range("C48").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D48").Value = "'9.53"
$ws.Range("E48").Value = "  +3.63%  "

$ws.Range("B49").Value = "Stacks"
$ws.Range("C49").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D49").Value = "'1.65"
$ws.Range("E49").Value = "  +9.20%  "

$ws.Range("B50").Value = "MultiversX"
$ws.Range("C50").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$ws.Range("D50").Value = "'54.24"
$ws.Range("E50").Value = "  +4.79%  "

$ws.Range("B51").Value = "BitcoinSV"
$ws.Range("C51").Value = "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
$ws.Range("D51").Value = "'74.06"
$ws.Range("E51").Value = "  +3.68%  "
